$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 873.0625
$ws.Range("I28").Value = 349
$ws.Range("K28").Value = 349
$ws.Range("M28").Value = 136
$ws.Range("H92").Value = 763.64703
$ws.Range("I92").Value = 784.3333
$ws.Range("K92").Value = 784.3333
$ws.Range("M92").Value = 463.6667
$ws.Range("H132").Value = 12794.319
$ws.Range("I132").Value = 2256.7646
$ws.Range("J132").Value = 22222.658
$ws.Range("K132").Value = 6770.293799999999
$ws.Range("L132").Value = 66667.974
$ws.Range("M132").Value = -4240.293799999999
$ws.Range("N132").Value = -71727.974
$ws.Range("H137").Value = 282418.25
$ws.Range("I137").Value = 834985.8
$ws.Range("J137").Value = 6134.4585
$ws.Range("K137").Value = 2504957.4
$ws.Range("L137").Value = 18403.3755
$ws.Range("M137").Value = -2502407.4
$ws.Range("N137").Value = -23503.3755
$ws.Range("H138").Value = 3408.8948
$ws.Range("I138").Value = 2644
$ws.Range("K138").Value = 7932
$ws.Range("M138").Value = -2792

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3660.1614
$ws.Range("I32").Value = 2025.8654
$ws.Range("J32").Value = 12158.5
$ws.Range("K32").Value = 2025.8654
$ws.Range("L32").Value = 12158.5
$ws.Range("M32").Value = -1738.8654
$ws.Range("N32").Value = -12732.5
$ws.Range("H45").Value = 4999
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("H61").Value = 14192.833
$ws.Range("I61").Value = 14192.833
$ws.Range("K61").Value = 14192.833
$ws.Range("M61").Value = -13980.833
$ws.Range("H74").Value = 20835582
$ws.Range("I74").Value = 31251562
$ws.Range("K74").Value = 31251562
$ws.Range("M74").Value = -31250688
$ws.Range("H77").Value = 20835582
$ws.Range("I77").Value = 31251562
$ws.Range("K77").Value = 156257810
$ws.Range("M77").Value = -156253442
$ws.Range("H110").Value = 5253.9
$ws.Range("I110").Value = 3840.3333
$ws.Range("J110").Value = 7374.25
$ws.Range("K110").Value = 3840.3333
$ws.Range("L110").Value = 7374.25
$ws.Range("M110").Value = -1795.3333
$ws.Range("N110").Value = -11464.25
$ws.Range("H122").Value = 5509.9
$ws.Range("I122").Value = 3633
$ws.Range("K122").Value = 10899
$ws.Range("M122").Value = -8449
$ws.Range("H136").Value = 14192.833
$ws.Range("I136").Value = 14192.833
$ws.Range("K136").Value = 42578.499
$ws.Range("M136").Value = -40028.499

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 4167192.5
$ws.Range("I64").Value = 6944891.5
$ws.Range("K64").Value = 6944891.5
$ws.Range("M64").Value = -6944666.5
$ws.Range("H67").Value = 4167192.5
$ws.Range("I67").Value = 6944891.5
$ws.Range("K67").Value = 6944891.5
$ws.Range("M67").Value = -6944111.5
$ws.Range("H86").Value = 1277.3478
$ws.Range("I86").Value = 1292
$ws.Range("K86").Value = 1292
$ws.Range("M86").Value = -169
$ws.Range("H89").Value = 1277.3478
$ws.Range("I89").Value = 1292
$ws.Range("K89").Value = 6460
$ws.Range("M89").Value = -844
$ws.Range("H99").Value = 4567768.5
$ws.Range("H134").Value = 4418.375
$ws.Range("I134").Value = 1302.1111
$ws.Range("J134").Value = 8425
$ws.Range("K134").Value = 3906.3333
$ws.Range("L134").Value = 25275
$ws.Range("M134").Value = -1371.3333
$ws.Range("N134").Value = -30345

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 21745182
$ws.Range("I31").Value = 76924600
$ws.Range("K31").Value = 76924600
$ws.Range("M31").Value = -76924305
$ws.Range("H34").Value = 21745182
$ws.Range("I34").Value = 76924600
$ws.Range("K34").Value = 76924600
$ws.Range("M34").Value = -76924398
$ws.Range("H105").Value = 2036.5555
$ws.Range("J105").Value = 2333
$ws.Range("L105").Value = 2333
$ws.Range("N105").Value = -5827

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 109.333336
$ws.Range("I2").Value = 85.5
$ws.Range("J2").Value = 123.35294
$ws.Range("K2").Value = 513
$ws.Range("L2").Value = 740.1176400000001
$ws.Range("M2").Value = -400
$ws.Range("N2").Value = -966.1176400000001
$ws.Range("H68").Value = 1724
$ws.Range("J68").Value = 1724
$ws.Range("L68").Value = 5172
$ws.Range("N68").Value = -6794
$ws.Range("H71").Value = 1724
$ws.Range("J71").Value = 1724
$ws.Range("L71").Value = 15516
$ws.Range("N71").Value = -23628
$ws.Range("H81").Value = 4087.3333
$ws.Range("I81").Value = 2572
$ws.Range("J81").Value = 5299.6
$ws.Range("K81").Value = 7716
$ws.Range("L81").Value = 15898.8
$ws.Range("M81").Value = -6593
$ws.Range("N81").Value = -18144.8
$ws.Range("H84").Value = 4087.3333
$ws.Range("I84").Value = 2572
$ws.Range("J84").Value = 5299.6
$ws.Range("K84").Value = 23148
$ws.Range("L84").Value = 47696.4
$ws.Range("M84").Value = -17532
$ws.Range("N84").Value = -58928.4
$ws.Range("H88").Value = 13333.333
$ws.Range("J88").Value = 13333.333
$ws.Range("L88").Value = 39999.999
$ws.Range("N88").Value = -40855.999
$ws.Range("H91").Value = 13333.333
$ws.Range("J91").Value = 13333.333
$ws.Range("L91").Value = 39999.999
$ws.Range("N91").Value = -42963.999
$ws.Range("H124").Value = 13000
$ws.Range("I124").Value = 8000
$ws.Range("J124").Value = 16333.333
$ws.Range("K124").Value = 24000
$ws.Range("L124").Value = 48999.999
$ws.Range("M124").Value = -19090
$ws.Range("N124").Value = -58819.999
$ws.Range("H132").Value = 5598.3335
$ws.Range("I132").Value = 717.25
$ws.Range("K132").Value = 6455.25
$ws.Range("M132").Value = -3925.25
$ws.Range("H137").Value = 7867747.5
$ws.Range("I137").Value = 940.55554
$ws.Range("J137").Value = 22028000
$ws.Range("K137").Value = 2821.66662
$ws.Range("L137").Value = 66084000
$ws.Range("M137").Value = 2278.33338
$ws.Range("N137").Value = -66094200
$ws.Range("H139").Value = 1868.65
$ws.Range("I139").Value = 1868.65
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 5605.950000000001
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -465.9500000000007
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 20841208
$ws.Range("I102").Value = 27785494
$ws.Range("K102").Value = 27785494
$ws.Range("M102").Value = -27783872
$ws.Range("H132").Value = 99881.664
$ws.Range("I132").Value = 185991.27
$ws.Range("J132").Value = 5161.1
$ws.Range("K132").Value = 557973.8099999999
$ws.Range("L132").Value = 15483.3
$ws.Range("M132").Value = -555443.8099999999
$ws.Range("N132").Value = -20543.3

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1896924.1
$ws.Range("I68").Value = 2527343
$ws.Range("J68").Value = 5667.3335
$ws.Range("K68").Value = 2527343
$ws.Range("L68").Value = 5667.3335
$ws.Range("M68").Value = -2526594
$ws.Range("N68").Value = -7165.3335
$ws.Range("H71").Value = 1896924.1
$ws.Range("I71").Value = 2527343
$ws.Range("J71").Value = 5667.3335
$ws.Range("K71").Value = 12636715
$ws.Range("L71").Value = 28336.6675
$ws.Range("M71").Value = -12632971
$ws.Range("N71").Value = -35824.6675
$ws.Range("H122").Value = 39413844
$ws.Range("I122").Value = 76927110
$ws.Range("K122").Value = 230781330
$ws.Range("M122").Value = -230778880
$ws.Range("H136").Value = 4066.5715
$ws.Range("I136").Value = 2284
$ws.Range("J136").Value = 6821.4546
$ws.Range("K136").Value = 6852
$ws.Range("L136").Value = 20464.3638
$ws.Range("M136").Value = -4302
$ws.Range("N136").Value = -25564.3638

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3161.5962
$ws.Range("I122").Value = 2749.9512
$ws.Range("J122").Value = 4695.909
$ws.Range("K122").Value = 8249.8536
$ws.Range("L122").Value = 14087.727
$ws.Range("M122").Value = -5799.8536
$ws.Range("N122").Value = -18987.727
$ws.Range("H132").Value = 30872100
$ws.Range("I132").Value = 7938542.5
$ws.Range("J132").Value = 45466180
$ws.Range("K132").Value = 23815627.5
$ws.Range("L132").Value = 136398540
$ws.Range("M132").Value = -23813097.5
$ws.Range("N132").Value = -136403600
$ws.Range("H136").Value = 10719.981
$ws.Range("I136").Value = 6344.4546
$ws.Range("K136").Value = 19033.3638
$ws.Range("M136").Value = -16483.3638
$ws.Range("H138").Value = 100000
$ws.Range("J138").Value = 100000
$ws.Range("L138").Value = 100000
$ws.Range("N138").Value = -110280
